$wb = $excel.ActiveWorkbook

$wsShipment = $wb.Worksheets.Item("ShipmentInformation")
$wsShipment.Range("C2").Value = "PickUp115"
$wsShipment.Range("K2").Value = "DropOff95"

$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B2").Value = "02-15-2022"
$wsInput.Range("T2").Value = "58572102"
$wsInput.Range("U2").Value = "`$709.82"
$wsInput.Range("W2").Value = "FCT943034220790415360"
$wsInput.Range("X2").Value = "FCTEST1004293"
